$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colE = @(
  "codeforiati:group-code",
  "110",
  "110",
  "110",
  "110",
  "110",
  "110",
  "110",
  "110",
  "110",
  "110",
  "110",
  "110",
  "110",
  "110",
  "110",
  "120",
  "120",
  "120",
  "120",
  "120",
  "120",
  "120",
  "120",
  "120",
  "120",
  "120",
  "120",
  "120",
  "120",
  "120",
  "120",
  "120",
  "120",
  "120",
  "130",
  "130",
  "130",
  "130",
  "130",
  "140",
  "140",
  "140",
  "140",
  "140",
  "140",
  "140",
  "140",
  "140",
  "140",
  "140",
  "150",
  "150",
  "150",
  "150",
  "150",
  "150",
  "150",
  "150",
  "150",
  "150",
  "150",
  "150",
  "150",
  "150",
  "150",
  "150",
  "150",
  "150",
  "150",
  "150",
  "150",
  "150",
  "160",
  "160",
  "160",
  "160",
  "160",
  "160",
  "160",
  "160",
  "160",
  "160",
  "160",
  "210",
  "210",
  "210",
  "210",
  "210",
  "210",
  "210",
  "220",
  "220",
  "220",
  "220",
  "230",
  "230",
  "230",
  "230",
  "230",
  "230",
  "230",
  "230",
  "230",
  "230",
  "230",
  "230",
  "230",
  "230",
  "230",
  "230",
  "230",
  "230",
  "230",
  "230",
  "230",
  "230",
  "230",
  "230",
  "230",
  "230",
  "230",
  "230",
  "240",
  "240",
  "240",
  "240",
  "240",
  "240",
  "250",
  "250",
  "250",
  "250",
  "310",
  "310",
  "310",
  "310",
  "310",
  "310",
  "310",
  "310",
  "310",
  "310",
  "310",
  "310",
  "310",
  "310",
  "310",
  "310",
  "310",
  "310",
  "310",
  "310",
  "310",
  "310",
  "310",
  "310",
  "310",
  "310",
  "310",
  "310",
  "310",
  "320",
  "320",
  "320",
  "320",
  "320",
  "320",
  "320",
  "320",
  "320",
  "320",
  "320",
  "320",
  "320",
  "320",
  "320",
  "320",
  "320",
  "320",
  "320",
  "320",
  "320",
  "320",
  "320",
  "320",
  "320",
  "320",
  "320",
  "320",
  "320",
  "320",
  "331",
  "331",
  "331",
  "331",
  "331",
  "331",
  "332",
  "410",
  "410",
  "410",
  "410",
  "410",
  "410",
  "430",
  "430",
  "430",
  "430",
  "430",
  "430",
  "430",
  "430",
  "430",
  "430",
  "510",
  "520",
  "530",
  "530",
  "600",
  "600",
  "600",
  "600",
  "600",
  "600",
  "600",
  "720",
  "720",
  "720",
  "730",
  "740",
  "910",
  "930",
  "998",
  "998"
)

$colF = @(
  "codeforiati:group-name",
  "Educación",
  "Educación",
  "Educación",
  "Educación",
  "Educación",
  "Educación",
  "Educación",
  "Educación",
  "Educación",
  "Educación",
  "Educación",
  "Educación",
  "Educación",
  "Educación",
  "Educación",
  "Salud",
  "Salud",
  "Salud",
  "Salud",
  "Salud",
  "Salud",
  "Salud",
  "Salud",
  "Salud",
  "Salud",
  "Salud",
  "Salud",
  "Salud",
  "Salud",
  "Salud",
  "Salud",
  "Salud",
  "Salud",
  "Salud",
  "Programas/políticas sobre población y salud reproductiva",
  "Programas/políticas sobre población y salud reproductiva",
  "Programas/políticas sobre población y salud reproductiva",
  "Programas/políticas sobre población y salud reproductiva",
  "Programas/políticas sobre población y salud reproductiva",
  "Abastecimiento de agua y saneamiento",
  "Abastecimiento de agua y saneamiento",
  "Abastecimiento de agua y saneamiento",
  "Abastecimiento de agua y saneamiento",
  "Abastecimiento de agua y saneamiento",
  "Abastecimiento de agua y saneamiento",
  "Abastecimiento de agua y saneamiento",
  "Abastecimiento de agua y saneamiento",
  "Abastecimiento de agua y saneamiento",
  "Abastecimiento de agua y saneamiento",
  "Abastecimiento de agua y saneamiento",
  "Gobierno y sociedad civil",
  "Gobierno y sociedad civil",
  "Gobierno y sociedad civil",
  "Gobierno y sociedad civil",
  "Gobierno y sociedad civil",
  "Gobierno y sociedad civil",
  "Gobierno y sociedad civil",
  "Gobierno y sociedad civil",
  "Gobierno y sociedad civil",
  "Gobierno y sociedad civil",
  "Gobierno y sociedad civil",
  "Gobierno y sociedad civil",
  "Gobierno y sociedad civil",
  "Gobierno y sociedad civil",
  "Gobierno y sociedad civil",
  "Gobierno y sociedad civil",
  "Gobierno y sociedad civil",
  "Gobierno y sociedad civil",
  "Gobierno y sociedad civil",
  "Gobierno y sociedad civil",
  "Gobierno y sociedad civil",
  "Gobierno y sociedad civil",
  "Otros servicios e infraestructuras sociales",
  "Otros servicios e infraestructuras sociales",
  "Otros servicios e infraestructuras sociales",
  "Otros servicios e infraestructuras sociales",
  "Otros servicios e infraestructuras sociales",
  "Otros servicios e infraestructuras sociales",
  "Otros servicios e infraestructuras sociales",
  "Otros servicios e infraestructuras sociales",
  "Otros servicios e infraestructuras sociales",
  "Otros servicios e infraestructuras sociales",
  "Otros servicios e infraestructuras sociales",
  "Transporte y almacenamiento",
  "Transporte y almacenamiento",
  "Transporte y almacenamiento",
  "Transporte y almacenamiento",
  "Transporte y almacenamiento",
  "Transporte y almacenamiento",
  "Transporte y almacenamiento",
  "Comunicaciones",
  "Comunicaciones",
  "Comunicaciones",
  "Comunicaciones",
  "Energía",
  "Energía",
  "Energía",
  "Energía",
  "Energía",
  "Energía",
  "Energía",
  "Energía",
  "Energía",
  "Energía",
  "Energía",
  "Energía",
  "Energía",
  "Energía",
  "Energía",
  "Energía",
  "Energía",
  "Energía",
  "Energía",
  "Energía",
  "Energía",
  "Energía",
  "Energía",
  "Energía",
  "Energía",
  "Energía",
  "Energía",
  "Energía",
  "Servicios bancarios y financieros",
  "Servicios bancarios y financieros",
  "Servicios bancarios y financieros",
  "Servicios bancarios y financieros",
  "Servicios bancarios y financieros",
  "Servicios bancarios y financieros",
  "Empresas y otros servicios",
  "Empresas y otros servicios",
  "Empresas y otros servicios",
  "Empresas y otros servicios",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Agricultura, Silvicultura, Pesca",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Industria, extractivas, construcción",
  "Política y regulación comercial",
  "Política y regulación comercial",
  "Política y regulación comercial",
  "Política y regulación comercial",
  "Política y regulación comercial",
  "Política y regulación comercial",
  "Turismo",
  "Protección general medio ambiente",
  "Protección general medio ambiente",
  "Protección general medio ambiente",
  "Protección general medio ambiente",
  "Protección general medio ambiente",
  "Protección general medio ambiente",
  "Otras acciones multisectoriales",
  "Otras acciones multisectoriales",
  "Otras acciones multisectoriales",
  "Otras acciones multisectoriales",
  "Otras acciones multisectoriales",
  "Otras acciones multisectoriales",
  "Otras acciones multisectoriales",
  "Otras acciones multisectoriales",
  "Otras acciones multisectoriales",
  "Otras acciones multisectoriales",
  "Apoyo presupuestario general",
  "Ayuda alimentaria para el desarrollo",
  "Otras ayudas en forma de suministro de bienes",
  "Otras ayudas en forma de suministro de bienes",
  "Actividades relacionadas con la deuda",
  "Actividades relacionadas con la deuda",
  "Actividades relacionadas con la deuda",
  "Actividades relacionadas con la deuda",
  "Actividades relacionadas con la deuda",
  "Actividades relacionadas con la deuda",
  "Actividades relacionadas con la deuda",
  "Ayuda de emergencia",
  "Ayuda de emergencia",
  "Ayuda de emergencia",
  "Ayuda a la reconstrucción y a la rehabilitación",
  "Prevención de desastres",
  "Costes administrativos donantes",
  "Ayuda a refugiados en el país donante",
  "Sin especificación / no clasificados",
  "Sin especificación / no clasificados"
)

$colG = @(
  "codeforiati:category-name",
  "Educación, nivel no especificado",
  "Educación, nivel no especificado",
  "Educación, nivel no especificado",
  "Educación, nivel no especificado",
  "Educación básica",
  "Educación básica",
  "Educación básica",
  "Educación básica",
  "Educación básica",
  "Educación básica",
  "Educación básica",
  "Educación secundaria",
  "Educación secundaria",
  "Educación post-secundaria",
  "Educación post-secundaria",
  "Salud, general",
  "Salud, general",
  "Salud, general",
  "Salud, general",
  "Salud básica",
  "Salud básica",
  "Salud básica",
  "Salud básica",
  "Salud básica",
  "Salud básica",
  "Salud básica",
  "Salud básica",
  "Salud básica",
  "Enfermedades no transmisibles",
  "Enfermedades no transmisibles",
  "Enfermedades no transmisibles",
  "Enfermedades no transmisibles",
  "Enfermedades no transmisibles",
  "Enfermedades no transmisibles",
  "Programas/políticas sobre población y salud reproductiva",
  "Programas/políticas sobre población y salud reproductiva",
  "Programas/políticas sobre población y salud reproductiva",
  "Programas/políticas sobre población y salud reproductiva",
  "Programas/políticas sobre población y salud reproductiva",
  "Abastecimiento de agua y saneamiento",
  "Abastecimiento de agua y saneamiento",
  "Abastecimiento de agua y saneamiento",
  "Abastecimiento de agua y saneamiento",
  "Abastecimiento de agua y saneamiento",
  "Abastecimiento de agua y saneamiento",
  "Abastecimiento de agua y saneamiento",
  "Abastecimiento de agua y saneamiento",
  "Abastecimiento de agua y saneamiento",
  "Abastecimiento de agua y saneamiento",
  "Abastecimiento de agua y saneamiento",
  "Gobierno y sociedad civil, general",
  "Gobierno y sociedad civil, general",
  "Gobierno y sociedad civil, general",
  "Gobierno y sociedad civil, general",
  "Gobierno y sociedad civil, general",
  "Gobierno y sociedad civil, general",
  "Gobierno y sociedad civil, general",
  "Gobierno y sociedad civil, general",
  "Gobierno y sociedad civil, general",
  "Gobierno y sociedad civil, general",
  "Gobierno y sociedad civil, general",
  "Gobierno y sociedad civil, general",
  "Gobierno y sociedad civil, general",
  "Gobierno y sociedad civil, general",
  "Gobierno y sociedad civil, general",
  "Gobierno y sociedad civil, general",
  "Prevención y resolución de conflictos, paz y seguridad",
  "Prevención y resolución de conflictos, paz y seguridad",
  "Prevención y resolución de conflictos, paz y seguridad",
  "Prevención y resolución de conflictos, paz y seguridad",
  "Prevención y resolución de conflictos, paz y seguridad",
  "Prevención y resolución de conflictos, paz y seguridad",
  "Otros servicios e infraestructuras sociales",
  "Otros servicios e infraestructuras sociales",
  "Otros servicios e infraestructuras sociales",
  "Otros servicios e infraestructuras sociales",
  "Otros servicios e infraestructuras sociales",
  "Otros servicios e infraestructuras sociales",
  "Otros servicios e infraestructuras sociales",
  "Otros servicios e infraestructuras sociales",
  "Otros servicios e infraestructuras sociales",
  "Otros servicios e infraestructuras sociales",
  "Otros servicios e infraestructuras sociales",
  "Transporte y almacenamiento",
  "Transporte y almacenamiento",
  "Transporte y almacenamiento",
  "Transporte y almacenamiento",
  "Transporte y almacenamiento",
  "Transporte y almacenamiento",
  "Transporte y almacenamiento",
  "Comunicaciones",
  "Comunicaciones",
  "Comunicaciones",
  "Comunicaciones",
  "Política energética",
  "Política energética",
  "Política energética",
  "Política energética",
  "Generación de energía, fuentes renovables",
  "Generación de energía, fuentes renovables",
  "Generación de energía, fuentes renovables",
  "Generación de energía, fuentes renovables",
  "Generación de energía, fuentes renovables",
  "Generación de energía, fuentes renovables",
  "Generación de energía, fuentes renovables",
  "Generación de energía, fuentes renovables",
  "Generación de energía, fuentes renovables",
  "Generación de energía, fuentes no renovables",
  "Generación de energía, fuentes no renovables",
  "Generación de energía, fuentes no renovables",
  "Generación de energía, fuentes no renovables",
  "Generación de energía, fuentes no renovables",
  "Generación de energía, fuentes no renovables",
  "Centrales de energía híbrida",
  "Centrales de energía nuclear",
  "Distribución de la energía",
  "Distribución de la energía",
  "Distribución de la energía",
  "Distribución de la energía",
  "Distribución de la energía",
  "Distribución de la energía",
  "Distribución de la energía",
  "Servicios bancarios y financieros",
  "Servicios bancarios y financieros",
  "Servicios bancarios y financieros",
  "Servicios bancarios y financieros",
  "Servicios bancarios y financieros",
  "Servicios bancarios y financieros",
  "Empresas y otros servicios",
  "Empresas y otros servicios",
  "Empresas y otros servicios",
  "Empresas y otros servicios",
  "Agricultura",
  "Agricultura",
  "Agricultura",
  "Agricultura",
  "Agricultura",
  "Agricultura",
  "Agricultura",
  "Agricultura",
  "Agricultura",
  "Agricultura",
  "Agricultura",
  "Agricultura",
  "Agricultura",
  "Agricultura",
  "Agricultura",
  "Agricultura",
  "Agricultura",
  "Agricultura",
  "Silvicultura",
  "Silvicultura",
  "Silvicultura",
  "Silvicultura",
  "Silvicultura",
  "Silvicultura",
  "Pesca",
  "Pesca",
  "Pesca",
  "Pesca",
  "Pesca",
  "Industria",
  "Industria",
  "Industria",
  "Industria",
  "Industria",
  "Industria",
  "Industria",
  "Industria",
  "Industria",
  "Industria",
  "Industria",
  "Industria",
  "Industria",
  "Industria",
  "Industria",
  "Industria",
  "Industria",
  "Industria",
  "Industria",
  "Industrias extractivas",
  "Industrias extractivas",
  "Industrias extractivas",
  "Industrias extractivas",
  "Industrias extractivas",
  "Industrias extractivas",
  "Industrias extractivas",
  "Industrias extractivas",
  "Industrias extractivas",
  "Industrias extractivas",
  "Construcción",
  "Política y regulación comercial",
  "Política y regulación comercial",
  "Política y regulación comercial",
  "Política y regulación comercial",
  "Política y regulación comercial",
  "Política y regulación comercial",
  "Turismo",
  "Protección general medio ambiente",
  "Protección general medio ambiente",
  "Protección general medio ambiente",
  "Protección general medio ambiente",
  "Protección general medio ambiente",
  "Protección general medio ambiente",
  "Otras acciones multisectoriales",
  "Otras acciones multisectoriales",
  "Otras acciones multisectoriales",
  "Otras acciones multisectoriales",
  "Otras acciones multisectoriales",
  "Otras acciones multisectoriales",
  "Otras acciones multisectoriales",
  "Otras acciones multisectoriales",
  "Otras acciones multisectoriales",
  "Otras acciones multisectoriales",
  "Apoyo presupuestario general",
  "Ayuda alimentaria para el desarrollo",
  "Otras ayudas en forma de suministro de bienes",
  "Otras ayudas en forma de suministro de bienes",
  "Actividades relacionadas con la deuda",
  "Actividades relacionadas con la deuda",
  "Actividades relacionadas con la deuda",
  "Actividades relacionadas con la deuda",
  "Actividades relacionadas con la deuda",
  "Actividades relacionadas con la deuda",
  "Actividades relacionadas con la deuda",
  "Ayuda de emergencia",
  "Ayuda de emergencia",
  "Ayuda de emergencia",
  "Ayuda a la reconstrucción y a la rehabilitación",
  "Prevención de desastres",
  "Costes administrativos donantes",
  "Ayuda a refugiados en el país donante",
  "Sin especificación / no clasificados",
  "Sin especificación / no clasificados"
)

for ($i = 0; $i -lt 235; $i++) {
  $r = $i + 1
  $ws.Cells.Item($r, 5).Value = $colE[$i]
  $ws.Cells.Item($r, 6).Value = $colF[$i]
  $ws.Cells.Item($r, 7).Value = $colG[$i]
}
"done"
